$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 341, pushing the existing rows 341-393 down to 342-394.
$ws.Rows.Item(341).Insert()

# Populate the newly inserted row 341 with the new weekly data record.
$ws.Cells.Item(341, 1).Value2  = 4
$ws.Cells.Item(341, 2).Value   = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(341, 3).Value   = "Los Lagos"
$ws.Cells.Item(341, 4).Value2  = 44984
$ws.Cells.Item(341, 5).Value2  = 10
$ws.Cells.Item(341, 6).Value2  = 100112043
$ws.Cells.Item(341, 7).Value   = "Pepino ensalada"
$ws.Cells.Item(341, 8).Value   = "Sin especificar"
$ws.Cells.Item(341, 9).Value   = "Primera"
$ws.Cells.Item(341, 10).Value2 = 80
$ws.Cells.Item(341, 11).Value2 = 13000
$ws.Cells.Item(341, 12).Value2 = 13000
$ws.Cells.Item(341, 13).Value2 = 13000
$ws.Cells.Item(341, 14).Value  = "$/caja 60 unidades"
$ws.Cells.Item(341, 15).Value  = "Región de Arica y Parinacota"
$ws.Cells.Item(341, 16).Value2 = 217
$ws.Cells.Item(341, 17).Value2 = 60
$ws.Cells.Item(341, 18).Value  = "Hortaliza"
